# Generate Report for Handback
# Rewrites the two tracked source-file GUID/names and their derived
# handoff/handback package filenames + timestamps across all three
# sheets (Overview, zh-cn, de-de), and keeps the hyperlink "display"
# text in sync with the new cell text.

$wb = $excel.ActiveWorkbook

$newMd1 = "af09b27d-3e2f-4e2e-bc32-4bbb023ea23c.md"
$newMd2 = "ffff6233187c-65f0-4fb4-bc3e-f7d4370b4c66.md"

# Both rows' handoff/handback package files now resolve to the SAME
# new package (built from the new GUID), for each locale.
$newXlfZh = "af09b27d-3e2f-4e2e-bc32-4bbb023ea23c.6e12c5a9b899230eba1fbb53799aae92711faa9a.zh-cn.xlf"
$newXlfDe = "af09b27d-3e2f-4e2e-bc32-4bbb023ea23c.6e12c5a9b899230eba1fbb53799aae92711faa9a.de-de.xlf"

$newDtZhE = "2016-03-18 00:49:53"
$newDtZhH = "2016-03-18 00:50:12"
$newDtDeE = "2016-03-18 00:49:56"
$newDtDeH = "2016-03-18 00:50:18"

# NOTE: this engine's PowerShell does not bind named (-Param value)
# arguments correctly, so helper functions below use positional params.
function Set-CellAndHyperlink {
    param($Sheet, [string]$CellAddr, [string]$NewText)
    $Sheet.Range($CellAddr).Value = $NewText
    $colPart = $CellAddr -replace '[0-9]+$', ''
    $rowPart = $CellAddr -replace '^[A-Za-z]+', ''
    $target = '$' + $colPart + '$' + $rowPart
    foreach ($hl in $Sheet.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq $target) {
            $hl.TextToDisplay = $NewText
        }
    }
}

# ---- Sheet "Overview" ----
$wsOverview = $wb.Worksheets.Item("Overview")
Set-CellAndHyperlink $wsOverview "A2" $newMd1
Set-CellAndHyperlink $wsOverview "A3" $newMd2

# ---- Sheet "zh-cn" ----
$wsZh = $wb.Worksheets.Item("zh-cn")
Set-CellAndHyperlink $wsZh "A2" $newMd1
Set-CellAndHyperlink $wsZh "D2" $newXlfZh
$wsZh.Range("E2").Value = $newDtZhE
Set-CellAndHyperlink $wsZh "F2" $newMd1
Set-CellAndHyperlink $wsZh "G2" $newXlfZh
$wsZh.Range("H2").Value = $newDtZhH

Set-CellAndHyperlink $wsZh "A3" $newMd2
Set-CellAndHyperlink $wsZh "D3" $newXlfZh
$wsZh.Range("E3").Value = $newDtZhE
Set-CellAndHyperlink $wsZh "F3" $newMd2
Set-CellAndHyperlink $wsZh "G3" $newXlfZh
$wsZh.Range("H3").Value = $newDtZhH

# ---- Sheet "de-de" ----
$wsDe = $wb.Worksheets.Item("de-de")
Set-CellAndHyperlink $wsDe "A2" $newMd1
Set-CellAndHyperlink $wsDe "D2" $newXlfDe
$wsDe.Range("E2").Value = $newDtDeE
Set-CellAndHyperlink $wsDe "F2" $newMd1
Set-CellAndHyperlink $wsDe "G2" $newXlfDe
$wsDe.Range("H2").Value = $newDtDeH

Set-CellAndHyperlink $wsDe "A3" $newMd2
Set-CellAndHyperlink $wsDe "D3" $newXlfDe
$wsDe.Range("E3").Value = $newDtDeE
Set-CellAndHyperlink $wsDe "F3" $newMd2
Set-CellAndHyperlink $wsDe "G3" $newXlfDe
$wsDe.Range("H3").Value = $newDtDeH
